# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# --- OFF sheet: Home row (row 2) updates ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 441
$wsOff.Range("C2").Value = 315
$wsOff.Range("D2").Value = 86
$wsOff.Range("E2").Value = 39

# --- DEF sheet: Home row (row 2) updates ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 451
$wsDef.Range("C2").Value = 306
$wsDef.Range("D2").Value = 107
$wsDef.Range("E2").Value = 54
$wsDef.Range("F2").Value = 9
